$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Dwellings_buildings")

# Row 2
$ws.Range("B2").Value = "All other industry"
$ws.Range("C2").Value = 220
$ws.Range("D2").Value = 235.8234527687297

# Row 3
$ws.Range("B3").Value = "Warehouses and storage"
$ws.Range("C3").Value = 220
$ws.Range("D3").Value = 182.2763843648209

# Row 4
$ws.Range("B4").Value = "Manufacturing and light industry"
$ws.Range("C4").Value = 220
$ws.Range("D4").Value = 219.0561889250815

# Row 5
$ws.Range("B5").Value = "All other industry"
$ws.Range("C5").Value = 480
$ws.Range("D5").Value = 235.8234527687297

# Row 6
$ws.Range("B6").Value = "Warehouses and storage"
$ws.Range("C6").Value = 480
$ws.Range("D6").Value = 182.2763843648209

# Row 7
$ws.Range("B7").Value = "Manufacturing and light industry"
$ws.Range("C7").Value = 480
$ws.Range("D7").Value = 219.0561889250815
